$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values recomputed (commit: "regen save_data to use K instead of
# Strike#, regen std/mean, calc and write s_vals"). Write the new K values for
# each data row (rows 2-23) directly, matching the target diff.
$kValues = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 1
    8  = 1
    9  = 2
    10 = 3
    11 = 4
    12 = 0
    13 = 2
    14 = 3
    15 = 1
    16 = 2
    17 = 1
    18 = 1
    19 = 2
    20 = 0
    21 = 2
    22 = 1
    23 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
